# issue #5: stock data output to json file
#
# The "股票" (stock) sheet gains a new "property_category" column (with
# constant value "stock" for every data row), inserted between the
# existing "total" and "date" columns. The columns that used to sit at
# H/I/J (date / legislator_name / legislator_id) shift right to I/J/K.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new column before the current "date" column (H), shifting
# date/legislator_name/legislator_id one column to the right.
$ws.Range("H1").EntireColumn.Insert()

# New header cell.
$ws.Range("H1").Value = "property_category"

# Every stock row on this sheet is categorized as "stock".
$ws.Range("H2:H13").Value = "stock"
